$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update "Förändrad" (column C) date for existing rows 2-42 from 2025-03-06 to 2025-03-07
for ($r = 2; $r -le 42; $r++) {
    $ws.Cells.Item($r, 3).Value = 45723
}

# Row 42 gets an explicit (custom) row height marker in the saved file
$ws.Rows.Item(42).RowHeight = 15

# Add new row 43 with data for case A 8649-2025
$ws.Cells.Item(43, 1).Value = "A 8649-2025"
$ws.Cells.Item(43, 2).Value = 45710
$ws.Cells.Item(43, 3).Value = 45723
$ws.Cells.Item(43, 4).Value = "OKÄNT"
$ws.Cells.Item(43, 5).Value = "OKÄNT"
$ws.Cells.Item(43, 7).Value = 0.9
$ws.Cells.Item(43, 8).Value = 0
$ws.Cells.Item(43, 9).Value = 0
$ws.Cells.Item(43, 10).Value = 0
$ws.Cells.Item(43, 11).Value = 0
$ws.Cells.Item(43, 12).Value = 0
$ws.Cells.Item(43, 13).Value = 0
$ws.Cells.Item(43, 14).Value = 0
$ws.Cells.Item(43, 15).Value = 0
$ws.Cells.Item(43, 16).Value = 0
$ws.Cells.Item(43, 17).Value = 0

# Copy style from row 42 (B, C date style; R wrap-text style)
$ws.Cells.Item(42, 2).Copy() | Out-Null
$ws.Cells.Item(43, 2).PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Cells.Item(42, 3).Copy() | Out-Null
$ws.Cells.Item(43, 3).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(42, 18).Copy() | Out-Null
$ws.Cells.Item(43, 18).PasteSpecial(-4122) | Out-Null
